$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
$headerIndexes = @(1, 2)
$footerIndexes = @(1, 2)

foreach ($idx in $headerIndexes) {
    $hdr = $sec.Headers.Item($idx)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image1.jpg"
            }
        }
    }
}

foreach ($idx in $footerIndexes) {
    $ftr = $sec.Footers.Item($idx)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -like "*PearsonLogo.png") {
                # Direct property access on footer InlineShapes can raise a stale
                # handle error in this runtime; routing the write through the
                # Selection object re-resolves the shape and succeeds.
                $shp.Range.Select()
                $word.Selection.InlineShapes.Item(1).Name = "image2.png"
            }
        }
    }
}

Write-Output "Rename complete"
